$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 590 rows: relabel Course Num to the merged "590&790" label and bump a couple of caps ---
$ws.Range("A39").Value = "590&790"
$ws.Range("C39").Value = 30

$ws.Range("A42").Value = "590&790"
$ws.Range("C42").Value = 50

$ws.Range("A44").Value = "590&790"

# --- 790 rows: renumber Sec # and refresh Enroll Cap ---
$ws.Range("B50").Value = 173
$ws.Range("C50").Value = 30

$ws.Range("B51").Value = 183
$ws.Range("C51").Value = 40

$ws.Range("B52").Value = 185

$ws.Range("B53").Value = 186
$ws.Range("C53").Value = 25

$ws.Range("B54").Value = 188
$ws.Range("C54").Value = 30

$ws.Range("B55").Value = 189
$ws.Range("C55").Value = 40

# Those sections (790/187/25, 790/188/30, 790/189/40) now live in the rows above,
# so clear the old trailing copies.
$ws.Range("A56:C58").ClearContents()

# Drop three now-redundant blank rows, shifting everything below up.
$ws.Range("A59:C61").EntireRow.Delete()

# Keep the AutoFilter / _FilterDatabase name / sort range in sync with the new extent.
$ws.AutoFilterMode = $false
$ws.Range("A1:C69").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='2025Fall'!`$A`$1:`$C`$69"
    }
}

# Sheet view: scroll position / active selection
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("E54").Select() | Out-Null
